# Fixed auto button problem
#
# Adds a new attendance sheet "2021-01-01" at the end of the workbook.
# The sheet reuses the exact same header/column layout and cell
# formatting (bold+bordered "Sr. No" header row, bold+bordered
# "Sr. No" data column) that every other daily attendance sheet in
# this workbook already uses.
#
# To guarantee the new sheet picks up the identical style definitions
# already present in the workbook (rather than Excel synthesizing new
# ones), we duplicate an existing same-shaped sheet (2020-11-21, which
# already has the same A1:I6 layout) and then simply overwrite its
# cell values with the new day's attendance data.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2020-11-21")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$templateSheet.Copy($null, $lastSheet)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2021-01-01"

# Header stays identical to the template, so only the data rows need
# to be replaced.
$data = @(
    @(4, "dishant", "802/ Gunjan nagar/ Andheri , Mumbai", "Team Lead", "18:35:31", 98.20766575855964, 147.3157019226992, "NA", "NA"),
    @(4, "dishant", "802/ Gunjan nagar/ Andheri , Mumbai", "Team Lead", "18:40:25", 96.85625977255769, 92.28159635919941, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "18:43:09", 95.5139015649184, 130.9243342687308, "NA", "NA"),
    @(1, "sachin", "301/Sanskruti-1,Andheri, Mumbai", "Software Engineer", "18:53:17", 97.63219807251328, 156.8040310916165, "NA", "NA"),
    @(4, "dishant", "802/ Gunjan nagar/ Andheri , Mumbai", "Team Lead", "18:54:01", 97.80702465077516, 109.9791580367084, "NA", "NA")
)

$r = 2
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}

$ws.Range("A1").Select()
